$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix cell formatting (style) drift first, by copying formats from cells
# that already carry the desired style index, before writing new values ---

# C7, C8, C9 need to switch from the "date" numFmt style to the plain
# center/wrap style already used by C2 (s=2).
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C7:C9").PasteSpecial(-4122) | Out-Null

# B11 needs to switch from its lone style (s=6) to the shared plain
# center/wrap style already used by B2 (s=2).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Row 7: Swap Salary (Update) ---
$ws.Range("E7").Value = "https://leetcode.com/problems/swap-salary/"
$ws.Range("D7").Value = "Swap Salary"
$ws.Range("B7").Value = "Update"
$ws.Range("C7").Value = "Easy"
$ws.Range("F7").Value = 45550
$ws.Range("G7").Value = "Sept"

# --- Row 8: Duplicate Emails (Select) ---
$ws.Range("E8").Value = "https://leetcode.com/problems/duplicate-emails/"
$ws.Range("D8").Value = "Duplicate-Emails"
$ws.Range("B8").Value = "Select"
$ws.Range("C8").Value = "Easy"
$ws.Range("F8").Value = 45550
$ws.Range("G8").Value = "Sept"

# --- Row 9: Employees Earning More Than Their Managers (Join) ---
$ws.Range("E9").Value = "https://leetcode.com/problems/employees-earning-more-than-their-managers/description/"
$ws.Range("D9").Value = " Employees Earning More Than Their Managers"
$ws.Range("B9").Value = "Join"
$ws.Range("C9").Value = "Easy"
$ws.Range("F9").Value = 45550
$ws.Range("G9").Value = "Sept"

# --- Row 10: Not Boring Movies (Select) ---
$ws.Range("E10").Value = "https://leetcode.com/problems/not-boring-movies/"
$ws.Range("D10").Value = "Not Boring Movies"
$ws.Range("B10").Value = "Select"
$ws.Range("C10").Value = "Easy"
$ws.Range("F10").Value = 45550
$ws.Range("G10").Value = "Sept"

# --- Row 11: Classes More Than 5 Students (Select) ---
$ws.Range("E11").Value = "https://leetcode.com/problems/classes-more-than-5-students/description/"
$ws.Range("D11").Value = "Classes More Than 5 Students"
$ws.Range("B11").Value = "Select"
$ws.Range("C11").Value = "Easy"
$ws.Range("F11").Value = 45550
$ws.Range("G11").Value = "Sept"

# --- Window / view cosmetics ---
$excel.ActiveWindow.Zoom = 112
$ws.Range("E13").Select() | Out-Null
